$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B12").Value = "ActivatedItem"
$ws.Range("B13").Value = "ActivatedItemMGr"

$ws.Range("B14").Select()
